# Gen1 Explicit – "Parents of Gen1 Retro writes to DB"
# Rearranges the Gen1LivedWIthFatherAtAgeX sheet: inserts a new variable
# record (R2837200) and re-sequences rows 2:41 so that the "326" block
# (R2837200-R2839100) comes first and the "306" block (R2839200-R2841100)
# comes second. Also flips which sheet/tab is active.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Gen1LivedWIthFatherAtAgeX  ")
$ws3 = $wb.Worksheets.Item("ParentHealth")

# Final (row, VariableCode, value) layout for B2:C41 after the edit.
$rowData = @(
  @(2, "R2837200", 326),
  @(3, "R2837300", 326),
  @(4, "R2837400", 326),
  @(5, "R2837500", 326),
  @(6, "R2837600", 326),
  @(7, "R2837700", 326),
  @(8, "R2837800", 326),
  @(9, "R2837900", 326),
  @(10, "R2838000", 326),
  @(11, "R2838100", 326),
  @(12, "R2838200", 326),
  @(13, "R2838300", 326),
  @(14, "R2838400", 326),
  @(15, "R2838500", 326),
  @(16, "R2838600", 326),
  @(17, "R2838700", 326),
  @(18, "R2838800", 326),
  @(19, "R2838900", 326),
  @(20, "R2839000", 326),
  @(21, "R2839100", 326),
  @(22, "R2839200", 306),
  @(23, "R2839300", 306),
  @(24, "R2839400", 306),
  @(25, "R2839500", 306),
  @(26, "R2839600", 306),
  @(27, "R2839700", 306),
  @(28, "R2839800", 306),
  @(29, "R2839900", 306),
  @(30, "R2840000", 306),
  @(31, "R2840100", 306),
  @(32, "R2840200", 306),
  @(33, "R2840300", 306),
  @(34, "R2840400", 306),
  @(35, "R2840500", 306),
  @(36, "R2840600", 306),
  @(37, "R2840700", 306),
  @(38, "R2840800", 306),
  @(39, "R2840900", 306),
  @(40, "R2841000", 306),
  @(41, "R2841100", 306),
)

foreach ($item in $rowData) {
  $r = $item[0]
  $code = $item[1]
  $val = $item[2]
  $ws1.Range("B$r").Value = $code
  $ws1.Range("C$r").Value = $val
}

# Row 2 is a brand-new record: it never had the old A2/J2 placeholder
# cells, so drop them entirely (rather than just clearing their content).
$ws1.Range("A2").Clear()
$ws1.Range("J2").Clear()

# Row 22 (the shifted-down first row of the old "306" block) now needs
# the same blank, formatted A/J placeholder cells row 2 used to have.
$ws1.Range("C2").Copy()
$ws1.Range("A22").PasteSpecial(-4122)
$ws1.Range("J22").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Gen1LivedWIthFatherAtAgeX becomes the active sheet/tab with A2:XFD41
# selected; ParentHealth reverts to its default (non-tab-selected) view.
$ws1.Activate()
$ws1.Range("A2:XFD41").Select()
